$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.589.28"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.233.43"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "269.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.01%  "
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0924"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.47%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.570.32"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "2.234.48"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "43.559.06"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +11.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.87%  "
$ws.Range("E28").Value = "  +5.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0930"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0350"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.09%  "
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("E39").Value = "  +19.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.217"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.438"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "2.453.96"
$ws.Range("E51").Value = "  +0.27%  "
